$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Cells whose new text would otherwise be auto-parsed as a plain number by Excel
# (no thousands separators) are written with a leading apostrophe so they stay text,
# exactly like the source data (prices are stored as text, not numbers).

$ws.Range('D2').Value = '61.538.61'
$ws.Range('E2').Value = '  +3.57%  '
$ws.Range('D3').Value = '3.062.86'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''574.12'
$ws.Range('E5').Value = '  +1.53%  '
$ws.Range('D6').Value = '''141.17'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.053.65'
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  +9.94%  '
$ws.Range('D12').Value = '''0.465'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = '''0.0000237'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').Value = '''34.78'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '3.570.27'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '''7.21'
$ws.Range('E17').Value = '  +2.48%  '
$ws.Range('D18').Value = '3.057.23'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').Value = '61.480.97'
$ws.Range('E19').Value = '  +3.68%  '
$ws.Range('D20').Value = '''447.64'
$ws.Range('E20').Value = '  +4.00%  '
$ws.Range('D21').Value = '''13.87'
$ws.Range('E21').Value = '  +1.87%  '
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Value = '''7.38'
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('D24').Value = '''13.45'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').Value = '''81.67'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('B27').Value = 'FirstDigitalUSD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '''2.21'
$ws.Range('E28').Value = '  +3.33%  '
$ws.Range('D29').Value = '''2.62'
$ws.Range('E29').Value = '  +2.67%  '
$ws.Range('D30').Value = '''7.98'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('D31').Value = '''6.64'
$ws.Range('E31').Value = '  +7.13%  '
$ws.Range('D32').Value = '''26.43'
$ws.Range('E33').Value = '  +8.25%  '
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('D35').Value = '0.0₃0788'
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('D36').Value = '''6.03'
$ws.Range('E36').Value = '  +4.06%  '
$ws.Range('D37').Value = '''2.15'
$ws.Range('E37').Value = '  +2.27%  '
$ws.Range('D38').Value = '''49.99'
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('D39').Value = '''2.95'
$ws.Range('E39').Value = '  +5.66%  '
$ws.Range('D40').Value = '''8.79'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').Value = '''419.91'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('E42').Value = '  +3.76%  '
$ws.Range('D43').Value = '2.763.25'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E45').Value = '  +6.31%  '
$ws.Range('D46').Value = '''35.63'
$ws.Range('E46').Value = '  +7.60%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''124.75'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('B49').Value = 'Fetch.AI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D49').Value = '''2.08'
$ws.Range('E49').Value = '  +1.97%  '
$ws.Range('E50').Value = '  +0.52%  '
$ws.Range('E51').Value = '  +0.28%  '
